# Apply the "MeanPDiffthdRanalysis" edit:
#  - rename the existing "Threshold" header (column E) to "oldThreshold"
#  - add two new columns: G ("Threshold") and H ("noobject") with data for
#    rows 2-13, and difference formulas (G = G-12row - H-12row) for rows 14-25

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row ---
$ws.Range("E1").Value = "oldThreshold"
$ws.Range("G1").Value = "Threshold"
$ws.Range("H1").Value = "noobject"

# --- New data columns G (Threshold) and H (noobject) for rows 2-13 ---
$gVals = @(0.78, 0.9, 0.82, 0.98, 1.1299999999999999, 1.27, 1.57, 2.25, 1.1399999999999999, 1.79, 1.1299999999999999, 1.43)
$hVals = @(0.38, 0.38, 0.45000000000000007, 0.44999999999999996, 0.7400000000000001, 0.74, 0.79, 0.79, 0.53, 0.53, 0.86999999999999988, 0.86999999999999988)

for ($i = 0; $i -lt 12; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 7).Value = $gVals[$i]
    $ws.Cells.Item($row, 8).Value = $hVals[$i]
}

# --- Difference formulas in column G for rows 14-25 (G(n) = G(n-12) - H(n-12)) ---
for ($row = 14; $row -le 25; $row++) {
    $srcRow = $row - 12
    $ws.Cells.Item($row, 7).Formula = "=G$srcRow-H$srcRow"
}
